# The deck is being stripped down to a single title slide (all the
# "Introduction to Git" content slides are removed and slide 1 becomes
# an empty "Title Slide" with ctrTitle / subTitle placeholders).

$p = $ppt.ActivePresentation

# 1. Remove every slide except the first one (slide 1 keeps its
#    identity / sldId so the presentation's sldIdLst ends up with just
#    the original <p:sldId id="256" r:id="rId2"/> entry).
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}

$s = $p.Slides.Item(1)

# 2. Switch slide 1 onto the "Title Slide" layout (index 1 in the
#    slide master's CustomLayouts, providing ctrTitle/subTitle
#    placeholders) instead of its current "Title and Content" layout.
$s.CustomLayout = $p.SlideMaster.CustomLayouts.Item(1)

# 3. The old "Title" / "Content Placeholder" shapes are still present
#    alongside the two new ones pulled in from the layout. Each of
#    them is a layout-required placeholder, so deleting it once just
#    re-creates an empty copy; the second delete removes it for good.
#    Repeat twice each to clear both leftover shapes and end up with
#    only the new ctrTitle + subTitle placeholders.
for ($j = 0; $j -lt 2; $j++) {
    $s.Shapes.Item(1).Delete()
    $s.Shapes.Item(1).Delete()
}
